$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.602.58"
$ws.Range("E2").Value = "'  -3.97%  "

# Row 3
$ws.Range("D3").Value = "'3.391.73"
$ws.Range("E3").Value = "'  -6.03%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'188.07"
$ws.Range("E5").Value = "'  -6.87%  "

# Row 6
$ws.Range("D6").Value = "'529.98"
$ws.Range("E6").Value = "'  -6.42%  "

# Row 7
$ws.Range("D7").Value = "'0.607"
$ws.Range("E7").Value = "'  -2.42%  "

# Row 8
$ws.Range("D8").Value = "'3.389.95"
$ws.Range("E8").Value = "'  -5.94%  "

# Row 9
$ws.Range("E9").Value = "'  -0.07%  "

# Row 10
$ws.Range("D10").Value = "'0.629"
$ws.Range("E10").Value = "'  -6.77%  "

# Row 11
$ws.Range("D11").Value = "'59.15"
$ws.Range("E11").Value = "'  -3.08%  "

# Row 12
$ws.Range("D12").Value = "'0.134"
$ws.Range("E12").Value = "'  -11.53%  "

# Row 13
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "'  -11.04%  "

# Row 14
$ws.Range("D14").Value = "'9.30"
$ws.Range("E14").Value = "'  -7.28%  "

# Row 15
$ws.Range("D15").Value = "'3.920.93"
$ws.Range("E15").Value = "'  -6.55%  "

# Row 16
$ws.Range("E16").Value = "'  -2.89%  "

# Row 17
$ws.Range("D17").Value = "'3.387.47"
$ws.Range("E17").Value = "'  -6.39%  "

# Row 18
$ws.Range("D18").Value = "'65.270.05"
$ws.Range("E18").Value = "'  -4.29%  "

# Row 19
$ws.Range("D19").Value = "'17.54"
$ws.Range("E19").Value = "'  -7.95%  "

# Row 20
$ws.Range("D20").Value = "'11.18"
$ws.Range("E20").Value = "'  -9.57%  "

# Row 21
$ws.Range("D21").Value = "'0.978"
$ws.Range("E21").Value = "'  -9.16%  "

# Row 22
$ws.Range("D22").Value = "'374.19"
$ws.Range("E22").Value = "'  -7.16%  "

# Row 23
$ws.Range("D23").Value = "'81.86"
$ws.Range("E23").Value = "'  -4.20%  "

# Row 24
$ws.Range("D24").Value = "'3.74"
$ws.Range("E24").Value = "'  -10.11%  "

# Row 25
$ws.Range("D25").Value = "'10.89"
$ws.Range("E25").Value = "'  -17.92%  "

# Row 26
$ws.Range("D26").Value = "'3.74"
$ws.Range("E26").Value = "'  -3.88%  "

# Row 27
$ws.Range("B27").Value = "'InternetComputer(DFINITY)"
$ws.Range("C27").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'11.61"
$ws.Range("E27").Value = "'  -7.86%  "

# Row 28
$ws.Range("B28").Value = "'ImmutableX"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.67"
$ws.Range("E28").Value = "'  -9.08%  "

# Row 29
$ws.Range("D29").Value = "'8.58"
$ws.Range("E29").Value = "'  -8.19%  "

# Row 30
$ws.Range("B30").Value = "'EthereumClassic"
$ws.Range("C30").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'29.79"
$ws.Range("E30").Value = "'  -5.63%  "

# Row 31
$ws.Range("B31").Value = "'Bittensor"
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'674.91"
$ws.Range("E31").Value = "'  +0.29%  "

# Row 32
$ws.Range("D32").Value = "'6.76"
$ws.Range("E32").Value = "'  -16.61%  "

# Row 33
$ws.Range("D33").Value = "'11.21"
$ws.Range("E33").Value = "'  -8.66%  "

# Row 34
$ws.Range("D34").Value = "'61.22"
$ws.Range("E34").Value = "'  -4.29%  "

# Row 35
$ws.Range("E35").Value = "'  -7.20%  "

# Row 36
$ws.Range("E36").Value = "'  +0.08%  "

# Row 37
$ws.Range("D37").Value = "'36.73"
$ws.Range("E37").Value = "'  -12.83%  "

# Row 38
$ws.Range("D38").Value = "'0.384"
$ws.Range("E38").Value = "'  -8.52%  "

# Row 39
$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "'  -0.34%  "

# Row 40
$ws.Range("E40").Value = "'  -6.25%  "

# Row 41
$ws.Range("D41").Value = "'2.869.62"
$ws.Range("E41").Value = "'  -11.42%  "

# Row 42
$ws.Range("D42").Value = "'2.78"
$ws.Range("E42").Value = "'  -13.04%  "

# Row 43
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "'  -6.88%  "

# Row 44
$ws.Range("D44").Value = "'0.0397"
$ws.Range("E44").Value = "'  -5.27%  "

# Row 45
$ws.Range("D45").Value = "'0.0₃0624"
$ws.Range("E45").Value = "'  -20.59%  "

# Row 46
$ws.Range("D46").Value = "'2.36"
$ws.Range("E46").Value = "'  -14.40%  "

# Row 47
$ws.Range("B47").Value = "'Monero"
$ws.Range("C47").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'138.24"
$ws.Range("E47").Value = "'  -0.67%  "

# Row 48
$ws.Range("B48").Value = "'Stellar"
$ws.Range("C48").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.126"
$ws.Range("E48").Value = "'  -4.48%  "

# Row 49
$ws.Range("B49").Value = "'ApeXProtocol"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "'  -6.30%  "

# Row 50
$ws.Range("B50").Value = "'Stacks"
$ws.Range("C50").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.63"
$ws.Range("E50").Value = "'  -4.32%  "

# Row 51
$ws.Range("D51").Value = "'7.74"
$ws.Range("E51").Value = "'  -12.86%  "
